# Add new columns I (I0) and J (IF) to the sheet, mirroring the
# existing header style (copied from H1) and filling in the per-row
# values for rows 2-76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from H1 onto
# the two new header cells before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for column I (and J, which is identical) for rows 2 through 76.
$values = @(7,6,5,7,7,6,8,8,6,7,3,8,7,8,7,8,7,8,7,8,7,6,7,8,8,9,9,9,8,8,9,8,9,8,8,9,9,9,9,9,9,8,8,9,8,9,8,11,10,7,9,9,9,8,8,9,8,9,9,8,8,8,9,9,9,9,8,6,6,6,5,7,3,7,5)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $val = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $val
    $ws.Cells.Item($row, 10).Value = $val
}
